# Auto-generated Excel COM-interop script
# Applies numeric corrections to the currentAveragePrice / Leve profit columns
# (columns H-N) on specific rows across the ALC, ARM, CRP, CUL, GSM, LTW, and WVR sheets,
# matching the values recorded in the upstream commit's OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Cells.Item(100, 8).Value = 2139.1538
$ws.Cells.Item(100, 9).Value = 1400.7142
$ws.Cells.Item(100, 11).Value = 1400.7142
$ws.Cells.Item(100, 13).Value = -859.7141999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 1623.9219
$ws.Cells.Item(132, 9).Value = 1368.9828
$ws.Cells.Item(132, 10).Value = 4088.3333
$ws.Cells.Item(132, 11).Value = 4106.9484
$ws.Cells.Item(132, 12).Value = 12264.9999
$ws.Cells.Item(132, 13).Value = -1576.9484
$ws.Cells.Item(132, 14).Value = -17324.9999
# Row 138
$ws.Cells.Item(138, 8).Value = 2302629.2
$ws.Cells.Item(138, 9).Value = 6668625.5
$ws.Cells.Item(138, 10).Value = 4736.4033
$ws.Cells.Item(138, 11).Value = 20005876.5
$ws.Cells.Item(138, 12).Value = 14209.2099
$ws.Cells.Item(138, 13).Value = -20000736.5
$ws.Cells.Item(138, 14).Value = -24489.2099

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 1799.8948
$ws.Cells.Item(61, 9).Value = 1470.75
$ws.Cells.Item(61, 11).Value = 1470.75
$ws.Cells.Item(61, 13).Value = -1258.75
# Row 63
$ws.Cells.Item(63, 8).Value = 5109.857
$ws.Cells.Item(63, 9).Value = 3826.5557
$ws.Cells.Item(63, 11).Value = 3826.5557
$ws.Cells.Item(63, 13).Value = -3140.5557
# Row 66
$ws.Cells.Item(66, 8).Value = 5109.857
$ws.Cells.Item(66, 9).Value = 3826.5557
$ws.Cells.Item(66, 11).Value = 19132.7785
$ws.Cells.Item(66, 13).Value = -15700.7785
# Row 136
$ws.Cells.Item(136, 8).Value = 1799.8948
$ws.Cells.Item(136, 9).Value = 1470.75
$ws.Cells.Item(136, 11).Value = 4412.25
$ws.Cells.Item(136, 13).Value = -1862.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2999.45
$ws.Cells.Item(31, 9).Value = 2311.5454
$ws.Cells.Item(31, 11).Value = 2311.5454
$ws.Cells.Item(31, 13).Value = -2016.5454
# Row 34
$ws.Cells.Item(34, 8).Value = 2999.45
$ws.Cells.Item(34, 9).Value = 2311.5454
$ws.Cells.Item(34, 11).Value = 2311.5454
$ws.Cells.Item(34, 13).Value = -2109.5454
# Row 58
$ws.Cells.Item(58, 8).Value = 1340.4688
$ws.Cells.Item(58, 9).Value = 1289.4828
$ws.Cells.Item(58, 11).Value = 1289.4828
$ws.Cells.Item(58, 13).Value = -1086.4828
# Row 86
$ws.Cells.Item(86, 8).Value = 2599.75
$ws.Cells.Item(86, 9).Value = 2100
$ws.Cells.Item(86, 10).Value = 3099.5
$ws.Cells.Item(86, 11).Value = 2100
$ws.Cells.Item(86, 12).Value = 3099.5
$ws.Cells.Item(86, 13).Value = -977
$ws.Cells.Item(86, 14).Value = -5345.5
# Row 89
$ws.Cells.Item(89, 8).Value = 2599.75
$ws.Cells.Item(89, 9).Value = 2100
$ws.Cells.Item(89, 10).Value = 3099.5
$ws.Cells.Item(89, 11).Value = 10500
$ws.Cells.Item(89, 12).Value = 15497.5
$ws.Cells.Item(89, 13).Value = -4884
$ws.Cells.Item(89, 14).Value = -26729.5
# Row 134
$ws.Cells.Item(134, 8).Value = 1571.2812
$ws.Cells.Item(134, 9).Value = 1355.8695
$ws.Cells.Item(134, 10).Value = 2121.7778
$ws.Cells.Item(134, 11).Value = 4067.6085
$ws.Cells.Item(134, 12).Value = 6365.3334
$ws.Cells.Item(134, 13).Value = -1532.6085
$ws.Cells.Item(134, 14).Value = -11435.3334
# Row 136
$ws.Cells.Item(136, 8).Value = 1340.4688
$ws.Cells.Item(136, 9).Value = 1289.4828
$ws.Cells.Item(136, 11).Value = 3868.4484
$ws.Cells.Item(136, 13).Value = -1318.4484

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1146.3914
$ws.Cells.Item(5, 9).Value = 1269.7273
$ws.Cells.Item(5, 10).Value = 1033.3334
$ws.Cells.Item(5, 11).Value = 3809.1819
$ws.Cells.Item(5, 12).Value = 3100.0002
$ws.Cells.Item(5, 13).Value = -3697.1819
$ws.Cells.Item(5, 14).Value = -3324.0002
# Row 12
$ws.Cells.Item(12, 8).Value = 840211.5600000001
$ws.Cells.Item(12, 10).Value = 1017075.5
$ws.Cells.Item(12, 12).Value = 3051226.5
$ws.Cells.Item(12, 14).Value = -3051572.5
# Row 15
$ws.Cells.Item(15, 8).Value = 205.2
$ws.Cells.Item(15, 9).Value = 75
$ws.Cells.Item(15, 10).Value = 400.5
$ws.Cells.Item(15, 11).Value = 225
$ws.Cells.Item(15, 12).Value = 1201.5
$ws.Cells.Item(15, 13).Value = -85
$ws.Cells.Item(15, 14).Value = -1481.5
# Row 22
$ws.Cells.Item(22, 8).Value = 23810692
$ws.Cells.Item(22, 9).Value = 83334080
$ws.Cells.Item(22, 10).Value = 1337.8
$ws.Cells.Item(22, 11).Value = 250002240
$ws.Cells.Item(22, 12).Value = 4013.4
$ws.Cells.Item(22, 13).Value = -250002071
$ws.Cells.Item(22, 14).Value = -4351.4
# Row 27
$ws.Cells.Item(27, 8).Value = 23810692
$ws.Cells.Item(27, 9).Value = 83334080
$ws.Cells.Item(27, 10).Value = 1337.8
$ws.Cells.Item(27, 11).Value = 250002240
$ws.Cells.Item(27, 12).Value = 4013.4
$ws.Cells.Item(27, 13).Value = -250002138
$ws.Cells.Item(27, 14).Value = -4217.4
# Row 113
$ws.Cells.Item(113, 8).Value = 916.9167
$ws.Cells.Item(113, 9).Value = 699.6667
$ws.Cells.Item(113, 11).Value = 2099.0001
$ws.Cells.Item(113, 13).Value = 70.9998999999998
# Row 117
$ws.Cells.Item(117, 8).Value = 64945.25
$ws.Cells.Item(117, 10).Value = 64945.25
$ws.Cells.Item(117, 12).Value = 194835.75
$ws.Cells.Item(117, 14).Value = -201719.75
# Row 122
$ws.Cells.Item(122, 8).Value = 721.2059
$ws.Cells.Item(122, 9).Value = 484.09525
$ws.Cells.Item(122, 10).Value = 1104.2307
$ws.Cells.Item(122, 11).Value = 4356.85725
$ws.Cells.Item(122, 12).Value = 9938.076300000001
$ws.Cells.Item(122, 13).Value = -1906.85725
$ws.Cells.Item(122, 14).Value = -14838.0763
# Row 131
$ws.Cells.Item(131, 8).Value = 889.49
$ws.Cells.Item(131, 10).Value = 889.49
$ws.Cells.Item(131, 12).Value = 2668.47
$ws.Cells.Item(131, 14).Value = -12748.47
# Row 132
$ws.Cells.Item(132, 8).Value = 1890.0312
$ws.Cells.Item(132, 9).Value = 1059.6154
$ws.Cells.Item(132, 10).Value = 2458.2104
$ws.Cells.Item(132, 11).Value = 9536.5386
$ws.Cells.Item(132, 12).Value = 22123.8936
$ws.Cells.Item(132, 13).Value = -7006.5386
$ws.Cells.Item(132, 14).Value = -27183.8936
# Row 135
$ws.Cells.Item(135, 8).Value = 1146.3914
$ws.Cells.Item(135, 9).Value = 1269.7273
$ws.Cells.Item(135, 10).Value = 1033.3334
$ws.Cells.Item(135, 11).Value = 11427.5457
$ws.Cells.Item(135, 12).Value = 9300.000599999999
$ws.Cells.Item(135, 13).Value = -8892.545700000001
$ws.Cells.Item(135, 14).Value = -14370.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 8333333.5
$ws.Cells.Item(11, 9).Value = 8333333.5
$ws.Cells.Item(11, 11).Value = 8333333.5
$ws.Cells.Item(11, 13).Value = -8333194.5
# Row 128
$ws.Cells.Item(128, 8).Value = 45740
$ws.Cells.Item(128, 10).Value = 45740
$ws.Cells.Item(128, 12).Value = 45740
$ws.Cells.Item(128, 14).Value = -55700
# Row 132
$ws.Cells.Item(132, 8).Value = 2260.2927
$ws.Cells.Item(132, 9).Value = 2039
$ws.Cells.Item(132, 10).Value = 2643.8667
$ws.Cells.Item(132, 11).Value = 6117
$ws.Cells.Item(132, 12).Value = 7931.6001
$ws.Cells.Item(132, 13).Value = -3587
$ws.Cells.Item(132, 14).Value = -12991.6001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3812.9285
$ws.Cells.Item(7, 9).Value = 3153.5557
$ws.Cells.Item(7, 10).Value = 4999.8
$ws.Cells.Item(7, 11).Value = 3153.5557
$ws.Cells.Item(7, 12).Value = 4999.8
$ws.Cells.Item(7, 13).Value = -3041.5557
$ws.Cells.Item(7, 14).Value = -5223.8
# Row 46
$ws.Cells.Item(46, 8).Value = 1415.7894
$ws.Cells.Item(46, 9).Value = 1016.6667
$ws.Cells.Item(46, 10).Value = 1775
$ws.Cells.Item(46, 11).Value = 1016.6667
$ws.Cells.Item(46, 12).Value = 1775
$ws.Cells.Item(46, 13).Value = -828.6667
$ws.Cells.Item(46, 14).Value = -2151
# Row 122
$ws.Cells.Item(122, 8).Value = 10421984
$ws.Cells.Item(122, 9).Value = 15630957
$ws.Cells.Item(122, 10).Value = 4038.125
$ws.Cells.Item(122, 11).Value = 46892871
$ws.Cells.Item(122, 12).Value = 12114.375
$ws.Cells.Item(122, 13).Value = -46890421
$ws.Cells.Item(122, 14).Value = -17014.375
# Row 126
$ws.Cells.Item(126, 8).Value = 3812.9285
$ws.Cells.Item(126, 9).Value = 3153.5557
$ws.Cells.Item(126, 10).Value = 4999.8
$ws.Cells.Item(126, 11).Value = 9460.667099999999
$ws.Cells.Item(126, 12).Value = 14999.4
$ws.Cells.Item(126, 13).Value = -6990.667099999999
$ws.Cells.Item(126, 14).Value = -19939.4

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Cells.Item(12, 8).Value = 27500
